$wb = $excel.ActiveWorkbook

# --- Sheet: rel_node__commodity (sheet7) ---
# Insert a new data row right after the header, shifting the existing
# 84 data rows down by one, and fill it with the new "313_HEAD_STORAGE" node.
$ws7 = $wb.Worksheets.Item("rel_node__commodity")
$ws7.Rows.Item(2).Insert() | Out-Null
$ws7.Range("A2").Value = "node__commodity"
$ws7.Range("B2").Value = "313_HEAD_STORAGE"
$ws7.Range("C2").Value = "electricity"

# --- Sheet: rel_node__temporal_block (sheet8) ---
# Two blocks of 84 rows each (blk_t1 then blk_t2). Insert a new row at the
# top of each block for the new node, shifting the remainder of that block
# (and everything after it) down by one row.
$ws8 = $wb.Worksheets.Item("rel_node__temporal_block")

# Top of the first block (rows 2-85, blk_t1)
$ws8.Rows.Item(2).Insert() | Out-Null
$ws8.Range("A2").Value = "node__temporal_block"
$ws8.Range("B2").Value = "313_HEAD_STORAGE"
$ws8.Range("C2").Value = "blk_t1"

# Top of the second block - after the first insertion the second block now
# starts at row 87 (was 86).
$ws8.Rows.Item(87).Insert() | Out-Null
$ws8.Range("A87").Value = "node__temporal_block"
$ws8.Range("B87").Value = "313_HEAD_STORAGE"
$ws8.Range("C87").Value = "blk_t2"

# --- View / selection state ---
# rel_node__commodity: selection sits on the newly inserted cell, sheet not active.
$ws7.Select()
$ws7.Range("B2").Select()

# rel_node__temporal_block: selection on the newly inserted second-block row,
# and this sheet ends up the active / tab-selected sheet.
$ws8.Select()
$ws8.Range("A87:C87").Select()

Write-Host "edit applied"
